$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.855.35"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +1.54%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.299.77"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +1.20%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.05%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'301.30"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +0.56%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'100.25"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +4.75%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.500"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +1.03%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  -0.06%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.510"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +4.10%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'35.85"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +8.27%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.0789"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +0.13%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  +2.38%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'  +11.42%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'6.82"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +2.26%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'2.665.97"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +1.47%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'2.315.14"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +1.97%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.800"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +1.78%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'42.754.07"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +1.40%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'12.45"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +6.31%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'6.17"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +3.14%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.0₃0897"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +0.69%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'67.68"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +2.05%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'235.42"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +0.03%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'2.21"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +12.93%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("B25").Value = "'Dai"
$ws.Range("B25").Style = "Normal"
$ws.Range("C25").Value = "'https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("C25").Style = "Normal"
$ws.Range("D25").Value = "'1.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -0.06%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("B26").Value = "'PancakeSwap"
$ws.Range("B26").Style = "Normal"
$ws.Range("C26").Value = "'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("C26").Style = "Normal"
$ws.Range("D26").Value = "'2.44"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -0.26%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'24.51"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +3.25%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'2.22"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +6.92%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("B29").Value = "'Monero"
$ws.Range("B29").Style = "Normal"
$ws.Range("C29").Value = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("C29").Style = "Normal"
$ws.Range("D29").Value = "'168.07"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +0.12%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("B30").Value = "'InjectiveProtocol"
$ws.Range("B30").Style = "Normal"
$ws.Range("C30").Value = "'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("C30").Style = "Normal"
$ws.Range("D30").Value = "'34.39"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +2.33%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'9.13"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +0.08%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  +0.11%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'4.97"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +1.45%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'4.62"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -1.41%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'17.25"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +3.27%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  +3.31%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'0.0686"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -0.33%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  +3.24%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'2.83"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +0.92%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'1.77"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +2.91%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.109"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +0.08%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'1.976.76"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +0.89%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.0286"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +3.33%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'2.20"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -4.33%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'10.15"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +5.88%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("B46").Value = "'EnergySwap"
$ws.Range("B46").Style = "Normal"
$ws.Range("C46").Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("C46").Style = "Normal"
$ws.Range("D46").Value = "'17.52"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -0.51%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("B47").Value = "'NEARProtocol"
$ws.Range("B47").Style = "Normal"
$ws.Range("C47").Value = "'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("C47").Style = "Normal"
$ws.Range("D47").Value = "'2.88"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +3.82%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'55.51"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +6.45%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'2.531.96"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +1.36%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'1.52"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +2.79%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("B51").Value = "'HuobiToken"
$ws.Range("B51").Style = "Normal"
$ws.Range("C51").Value = "'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("C51").Style = "Normal"
$ws.Range("D51").Value = "'2.73"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -0.29%  "
$ws.Range("E51").Style = "Normal"
